$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly record was added to the top of the Berenjena data block.
# Insert a fresh row at 184, which pushes the former rows 184:223 down
# to 185:224 (same shift the author's diff shows row-by-row).
$ws.Rows(184).Insert()

# Populate the newly inserted row with the new record's data.
$ws.Range("A184").Value = 6
$ws.Range("B184").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C184").Value = "Metropolitana"
$ws.Range("D184").Value = 44754
$ws.Range("E184").Value = 13
$ws.Range("F184").Value = 100112001
$ws.Range("G184").Value = "Berenjena"
$ws.Range("H184").Value = "Sin especificar"
$ws.Range("I184").Value = "Primera"
$ws.Range("J184").Value = 250
$ws.Range("K184").Value = 7000
$ws.Range("L184").Value = 8000
$ws.Range("M184").Value = 7400
$ws.Range("N184").Value = "`$/caja 50 unidades"
$ws.Range("O184").Value = "Región de Arica y Parinacota"
$ws.Range("P184").Value = 148
$ws.Range("Q184").Value = 50
$ws.Range("R184").Value = "Hortaliza"
